$d = $word.ActiveDocument

# Namespaces/wrapper used for every InsertXML call below: a pkg:package envelope
# that targets word/document.xml, carrying plain WordprocessingML body fragments.
function New-DocPackage($innerBodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerBodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# 1) Remove the "2009 Brighton Circle" address line entirely (run + its
#    paragraph-mark run properties / fr-FR language tag), leaving an empty
#    right-justified paragraph.
$rng = $d.Content
$rng.Find.Execute("2009 Brighton Circle")
$rng.InsertXML((New-DocPackage '<w:p><w:pPr><w:jc w:val="right"/></w:pPr></w:p>'))

# 2) "Augusta, GA 30906" -> "Augusta, GA " (drop the zip code, keep trailing space).
$d.Content.Find.Execute("Augusta, GA 30906", $true, $false, $false, $false, $false, $true, 1, $false, "Augusta, GA ", 2)

# 3) Split "Front-End Web Development Techdegree" into two runs, wrapping the
#    coined word "Techdegree" in spellcheck proofErr markers.
$rng = $d.Content
$rng.Find.Execute("Front-End Web Development Techdegree")
$rng.InsertXML((New-DocPackage ('<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Front-End Web Development </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Techdegree</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '</w:p>')))

# 4) Split "Associate Degree In Applied Science Major in Computer Programming"
#    into three runs, lower-casing "In" -> "in" as its own run.
$rng = $d.Content
$rng.Find.Execute("Associate Degree In Applied Science Major in Computer Programming")
$rng.InsertXML((New-DocPackage ('<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Associate Degree </w:t></w:r>' + `
    '<w:r><w:t>in</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> Applied Science Major in Computer Programming</w:t></w:r>' + `
    '</w:p>')))

# 5) Split " and Chrome DevTools" (second run of the bullet, after "Debugging in
#    Visual Studio") into " and Chrome " + "DevTools", wrapping "DevTools" in
#    spellcheck proofErr markers. Re-supply the whole paragraph because
#    InsertXML on a <w:p> replaces the full containing paragraph.
$rng = $d.Content
$rng.Find.Execute("Debugging in Visual Studio and Chrome DevTools")
$rng.InsertXML((New-DocPackage ('<w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr>' + `
    '<w:r><w:t>Debugging in Visual Studio</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> and Chrome </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>DevTools</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '</w:p>')))
